$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so numeric-looking
# strings (e.g. "1.009", "248.00") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 16 and 17 swapped rank order (Avalanche <-> BitcoinCash)
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# Updated Price (D) and Volume(1h) (E) figures for every coin row
$ws.Range("D2").Value = "30.458.77"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "1.939.91"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +1.29%  "
$ws.Range("D5").Value = "248.00"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "0.6992"
$ws.Range("E6").Value = "  -13.72%  "
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("D8").Value = "0.3247"
$ws.Range("E8").Value = "  -5.25%  "
$ws.Range("D9").Value = "26.55"
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").Value = "0.06788"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").Value = "0.7964"
$ws.Range("E11").Value = "  -6.60%  "
$ws.Range("D12").Value = "0.07996"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "1.958.86"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "5.389"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "94.13"
$ws.Range("E15").Value = "  -7.99%  "
$ws.Range("D16").Value = "263.14"
$ws.Range("E16").Value = "  -5.81%  "
$ws.Range("D17").Value = "14.53"
$ws.Range("E17").Value = "  +4.67%  "
$ws.Range("D18").Value = "30.505.88"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "5.880"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "0.000007816"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "2.220.02"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "1.008"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "6.847"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "9.667"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "158.70"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "2.261"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").Value = "0.1299"
$ws.Range("E29").Value = "  -20.55%  "
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "1.562"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "4.412"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").Value = "4.239"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").Value = "0.05103"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "1.193"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "0.7498"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "2.735"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "0.01940"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "2.792"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").Value = "80.14"
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").Value = "6.589"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "2.055"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").Value = "0.4432"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("D44").Value = "1.010"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "0.8421"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "101.95"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "9.745"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "7.313"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "36.09"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "1.494"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("D51").Value = "2.823"
$ws.Range("E51").Value = "  +32.08%  "
